$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.156.38'
$ws.Range('E2').Value = '  -3.29%  '
$ws.Range('D3').Value = '1.713.87'
$ws.Range('E3').Value = '  -3.66%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.11'
$ws.Range('E5').Value = '  -6.05%  '
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4768'
$ws.Range('E7').Value = '  +5.91%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3450'
$ws.Range('E8').Value = '  -3.32%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.13'
$ws.Range('E9').Value = '  +0.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07265'
$ws.Range('E10').Value = '  -2.55%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.044'
$ws.Range('E11').Value = '  -5.97%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.83'
$ws.Range('E13').Value = '  -5.73%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.853'
$ws.Range('E14').Value = '  -3.47%  '
$ws.Range('D15').Value = '1.716.95'
$ws.Range('E15').Value = '  -3.43%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.858'
$ws.Range('E16').Value = '  -5.74%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '88.79'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001038'
$ws.Range('E18').Value = '  -2.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06360'
$ws.Range('E19').Value = '  -1.35%  '
$ws.Range('E20').Value = '  +0.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.47'
$ws.Range('E21').Value = '  -3.99%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.609'
$ws.Range('E22').Value = '  -3.35%  '
$ws.Range('D23').Value = '27.190.83'
$ws.Range('E23').Value = '  -3.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.81'
$ws.Range('E24').Value = '  -4.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.089'
$ws.Range('E25').Value = '  -1.58%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.83'
$ws.Range('E26').Value = '  -6.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.62'
$ws.Range('E27').Value = '  -3.94%  '
$ws.Range('D28').Value = '1.912.26'
$ws.Range('E28').Value = '  -3.51%  '
$ws.Range('E29').Value = '  -3.60%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '119.87'
$ws.Range('E30').Value = '  -4.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.015'
$ws.Range('E31').Value = '  -8.43%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09264'
$ws.Range('E32').Value = '  +0.49%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.581'
$ws.Range('E33').Value = '  -2.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.305'
$ws.Range('E34').Value = '  -7.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02196'
$ws.Range('E35').Value = '  -4.38%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05904'
$ws.Range('E36').Value = '  -4.90%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '11.06'
$ws.Range('E37').Value = '  -7.20%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2003'
$ws.Range('E38').Value = '  -5.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.419'
$ws.Range('E39').Value = '  +1.53%  '
$ws.Range('B40').Value = 'Frax'
$ws.Range('C40').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9999'
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.731'
$ws.Range('E41').Value = '  -5.53%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5924'
$ws.Range('E42').Value = '  -6.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.094'
$ws.Range('E43').Value = '  -7.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.468'
$ws.Range('E44').Value = '  -5.85%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.65'
$ws.Range('E45').Value = '  -5.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.572'
$ws.Range('E46').Value = '  -4.84%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5607'
$ws.Range('E47').Value = '  -5.18%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '118.58'
$ws.Range('E48').Value = '  -3.50%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.835'
$ws.Range('E49').Value = '  -6.57%  '
$ws.Range('E50').Value = '  -3.79%  '
$ws.Range('E51').Value = '  -5.22%  '
